$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H33").Value = 215.75
$ws.Range("I33").Value = 121.111115
$ws.Range("K33").Value = 121.111115
$ws.Range("M33").Value = 107.888885
$ws.Range("H64").Value = 46146.39
$ws.Range("I64").Value = 202294
$ws.Range("K64").Value = 202294
$ws.Range("M64").Value = -202046
$ws.Range("H67").Value = 46146.39
$ws.Range("I67").Value = 202294
$ws.Range("K67").Value = 202294
$ws.Range("M67").Value = -201436
$ws.Range("H76").Value = 3496.8572
$ws.Range("I76").Value = 3496.8572
$ws.Range("K76").Value = 3496.8572
$ws.Range("M76").Value = -3181.8572
$ws.Range("H79").Value = 3496.8572
$ws.Range("I79").Value = 3496.8572
$ws.Range("K79").Value = 3496.8572
$ws.Range("M79").Value = -2404.8572
$ws.Range("H88").Value = 11100.25
$ws.Range("I88").Value = 5957.143
$ws.Range("J88").Value = 15100.444
$ws.Range("K88").Value = 5957.143
$ws.Range("L88").Value = 15100.444
$ws.Range("M88").Value = -5551.143
$ws.Range("N88").Value = -15912.444
$ws.Range("H91").Value = 11100.25
$ws.Range("I91").Value = 5957.143
$ws.Range("J91").Value = 15100.444
$ws.Range("K91").Value = 5957.143
$ws.Range("L91").Value = 15100.444
$ws.Range("M91").Value = -4553.143
$ws.Range("N91").Value = -17908.444
$ws.Range("H99").Value = 2925.3
$ws.Range("I99").Value = 2336.375
$ws.Range("J99").Value = 5281
$ws.Range("K99").Value = 7009.125
$ws.Range("L99").Value = 15843
$ws.Range("M99").Value = -5511.125
$ws.Range("N99").Value = -18839
$ws.Range("H121").Value = 1042367.1
$ws.Range("J121").Value = 1150170.9
$ws.Range("L121").Value = 3450512.7
$ws.Range("N121").Value = -3454006.7
$ws.Range("H126").Value = 45900.332
$ws.Range("J126").Value = 45900.332
$ws.Range("L126").Value = 45900.332
$ws.Range("N126").Value = -55780.332
$ws.Range("H128").Value = 47402.8
$ws.Range("J128").Value = 47402.8
$ws.Range("L128").Value = 47402.8
$ws.Range("N128").Value = -57362.8
$ws.Range("H130").Value = 51450.668
$ws.Range("J130").Value = 51450.668
$ws.Range("L130").Value = 51450.668
$ws.Range("N130").Value = -61490.668

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H63").Value = 2749.2856
$ws.Range("I63").Value = 2436.923
$ws.Range("J63").Value = 3020
$ws.Range("K63").Value = 2436.923
$ws.Range("L63").Value = 3020
$ws.Range("M63").Value = -1750.923
$ws.Range("N63").Value = -4392
$ws.Range("H66").Value = 2749.2856
$ws.Range("I66").Value = 2436.923
$ws.Range("J66").Value = 3020
$ws.Range("K66").Value = 12184.615
$ws.Range("L66").Value = 15100
$ws.Range("M66").Value = -8752.614999999998
$ws.Range("N66").Value = -21964
$ws.Range("H80").Value = 58641.332
$ws.Range("J80").Value = 58641.332
$ws.Range("L80").Value = 58641.332
$ws.Range("N80").Value = -60637.332
$ws.Range("H83").Value = 58641.332
$ws.Range("J83").Value = 58641.332
$ws.Range("L83").Value = 175923.996
$ws.Range("N83").Value = -185907.996
$ws.Range("H128").Value = 50421
$ws.Range("J128").Value = 50421
$ws.Range("L128").Value = 50421
$ws.Range("N128").Value = -60381
$ws.Range("H130").Value = 41295.8
$ws.Range("J130").Value = 41295.8
$ws.Range("L130").Value = 41295.8
$ws.Range("N130").Value = -51335.8

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H29").Value = 1000
$ws.Range("I29").Value = 1000
$ws.Range("K29").Value = 1000
$ws.Range("M29").Value = -711
$ws.Range("H35").Value = 33510
$ws.Range("J35").Value = 33510
$ws.Range("L35").Value = 33510
$ws.Range("N35").Value = -34130
$ws.Range("H82").Value = 11540.389
$ws.Range("I82").Value = 3080.7778
$ws.Range("J82").Value = 20000
$ws.Range("K82").Value = 3080.7778
$ws.Range("L82").Value = 20000
$ws.Range("M82").Value = -2697.7778
$ws.Range("N82").Value = -20766
$ws.Range("H85").Value = 11540.389
$ws.Range("I85").Value = 3080.7778
$ws.Range("J85").Value = 20000
$ws.Range("K85").Value = 3080.7778
$ws.Range("L85").Value = 20000
$ws.Range("M85").Value = -1754.7778
$ws.Range("N85").Value = -22652
$ws.Range("H86").Value = 2032.0769
$ws.Range("I86").Value = 2385.5715
$ws.Range("J86").Value = 1619.6666
$ws.Range("K86").Value = 2385.5715
$ws.Range("L86").Value = 1619.6666
$ws.Range("M86").Value = -1262.5715
$ws.Range("N86").Value = -3865.6666
$ws.Range("H89").Value = 2032.0769
$ws.Range("I89").Value = 2385.5715
$ws.Range("J89").Value = 1619.6666
$ws.Range("K89").Value = 11927.8575
$ws.Range("L89").Value = 8098.333000000001
$ws.Range("M89").Value = -6311.8575
$ws.Range("N89").Value = -19330.333
$ws.Range("H120").Value = 44761
$ws.Range("J120").Value = 44761
$ws.Range("L120").Value = 44761
$ws.Range("N120").Value = -54437
$ws.Range("H126").Value = 50772
$ws.Range("J126").Value = 50772
$ws.Range("L126").Value = 50772
$ws.Range("N126").Value = -60652
$ws.Range("H130").Value = 46063.57
$ws.Range("J130").Value = 46063.57
$ws.Range("L130").Value = 46063.57
$ws.Range("N130").Value = -56103.57

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H20").Value = 40283.168
$ws.Range("J20").Value = 40283.168
$ws.Range("L20").Value = 40283.168
$ws.Range("N20").Value = -40755.168
$ws.Range("H30").Value = 40283.168
$ws.Range("J30").Value = 40283.168
$ws.Range("L30").Value = 40283.168
$ws.Range("N30").Value = -40465.168
$ws.Range("I62").Value = 2700
$ws.Range("J62").Value = 5500
$ws.Range("K62").Value = 2700
$ws.Range("L62").Value = 5500
$ws.Range("M62").Value = -2076
$ws.Range("N62").Value = -6748
$ws.Range("I65").Value = 2700
$ws.Range("J65").Value = 5500
$ws.Range("K65").Value = 13500
$ws.Range("L65").Value = 27500
$ws.Range("M65").Value = -10380
$ws.Range("N65").Value = -33740
$ws.Range("H128").Value = 40283.168
$ws.Range("J128").Value = 40283.168
$ws.Range("L128").Value = 40283.168
$ws.Range("N128").Value = -50243.168

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H70").Value = 5447.619
$ws.Range("I70").Value = 5562.5
$ws.Range("J70").Value = 5080
$ws.Range("K70").Value = 5562.5
$ws.Range("L70").Value = 5080
$ws.Range("M70").Value = -5292.5
$ws.Range("N70").Value = -5620
$ws.Range("H73").Value = 5447.619
$ws.Range("I73").Value = 5562.5
$ws.Range("J73").Value = 5080
$ws.Range("K73").Value = 5562.5
$ws.Range("L73").Value = 5080
$ws.Range("M73").Value = -4626.5
$ws.Range("N73").Value = -6952
$ws.Range("H80").Value = 560977.75
$ws.Range("J80").Value = 8000
$ws.Range("L80").Value = 8000
$ws.Range("N80").Value = -9996
$ws.Range("H83").Value = 560977.75
$ws.Range("J83").Value = 8000
$ws.Range("L83").Value = 40000
$ws.Range("N83").Value = -49984
$ws.Range("H107").Value = 226294.89
$ws.Range("J107").Value = 12034.333
$ws.Range("L107").Value = 12034.333
$ws.Range("N107").Value = -15874.333
$ws.Range("H122").Value = 2093.75
$ws.Range("I122").Value = 2320
$ws.Range("J122").Value = 1716.6666
$ws.Range("K122").Value = 6960
$ws.Range("L122").Value = 5149.9998
$ws.Range("M122").Value = -4510
$ws.Range("N122").Value = -10049.9998
$ws.Range("H130").Value = 44183.5
$ws.Range("J130").Value = 44183.5
$ws.Range("L130").Value = 44183.5
$ws.Range("N130").Value = -54223.5

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H128").Value = 46214.5
$ws.Range("J128").Value = 46214.5
$ws.Range("L128").Value = 46214.5
$ws.Range("N128").Value = -56174.5

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H120").Value = 42476.8
$ws.Range("J120").Value = 42476.8
$ws.Range("L120").Value = 42476.8
$ws.Range("N120").Value = -52152.8
$ws.Range("H132").Value = 1213.5116
$ws.Range("I132").Value = 883.7941
$ws.Range("J132").Value = 2459.111
$ws.Range("K132").Value = 2651.3823
$ws.Range("L132").Value = 7377.333
$ws.Range("M132").Value = -121.3822999999998
$ws.Range("N132").Value = -12437.333
